$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear existing cell contents (keeps cell styles/formatting intact)
# so the shared-string table is rebuilt cleanly in the new order below.
$ws.Cells.ClearContents()

$names = @(
  'Cluster Name',
  'Active cases',
  '3035 Campbell Place Aged Care Glen Waverley',
  '3364 Assisi Centre Aged Care Rosanna',
  '3528 Ottoman Village Aged Care Broadmeadows',
  '3622 Olivet Care Aged Care Services Ringwood',
  '3633 Lifeview Emerald Glades Aged Care Emerald',
  '3652 Regis Aged Care Dandenong North',
  '3824 Estia Health South Morang',
  '3961 Heritage Care Water Gardens Aged Care Facility Sydenham',
  'AG Industries Pty Ltd Factory Thomastown',
  'Aintree Primary School Aintree',
  'Athol Road Primary School Springvale South',
  'Australian Meat Group Abattoir Dandenong South',
  'Bacchus Marsh Childcare and Kindergarten Centre Bacchus Marsh',
  'Baden Powell College Tarneit',
  'Bandiana Primary School Bandiana',
  'CREST Children''s Sanctuary Dandenong',
  'Dandenong South Primary School Dandenong',
  'Hamlyn Views School Hamlyn Heights',
  'Hazelwood North Primary School Hazelwood North',
  'Hippity Hop Childcare and Kindergarten Pakenham',
  'Life Church Mooroopna',
  'Lilydale Motor Inn Lilydale',
  'Lowanna College Newborough',
  'M.C. Herd Corio',
  'Master Poultry Group West Footscray',
  'Morwell Park Primary School Morwell',
  'Northern Bay College Goldsworthy 9-12 Campus Corio',
  'Northern Bay College Wexford Campus Corio',
  'Oakleigh South Primary School Oakleigh South',
  'Saint Augustines Primary School Wodonga',
  'Saint Monica''s Primary School Wodonga',
  'Smartie Pants Early Learning and Development Diamond Creek',
  'St Georges Road Primary School Shepparton',
  'St Josephs Catholic Primary School Warragul',
  'St Josephs Primary School Quarry Hill',
  'St Mary''s Primary School Swan Hill',
  'St Vincents Hospital Emergency Department Melbourne',
  'St. Brendan''s Catholic Primary School Lakes Entrance',
  'Stockdale Road Primary School Traralgon',
  'Sunbury Primary School Sunbury',
  'TUROSI PTY LTD Thomastown',
  'The Lake Primary School Cabarita',
  'Werribee Mercy Hospital Emergency Department',
  'Western Health Sunshine Hospital Emergency Department St Albans',
  'Wodonga Cemetery Wodonga',
  'Wodonga Middle Years College Huon Campus Wodonga',
  'Wodonga Primary School Wodonga',
  'Wodonga Senior Secondary College Wodonga',
  'Wodonga South Primary School Wodonga',
  'Woodend Primary School Woodend',
  'Yooralla Disability Residential Care Alfrieda Street St Albans'
)

$values = @{
  2 = 13
  3 = 20
  4 = 12
  5 = 14
  6 = 17
  7 = 12
  8 = 17
  9 = 23
  10 = 13
  11 = 18
  12 = 11
  13 = 21
  14 = 22
  15 = 12
  16 = 11
  17 = 10
  18 = 11
  19 = 13
  20 = 13
  21 = 10
  22 = 12
  23 = 10
  24 = 12
  25 = 10
  26 = 11
  27 = 25
  28 = 16
  29 = 55
  30 = 16
  31 = 10
  32 = 11
  33 = 11
  34 = 12
  35 = 12
  36 = 17
  37 = 10
  38 = 17
  39 = 10
  40 = 23
  41 = 10
  42 = 10
  43 = 12
  44 = 33
  45 = 14
  46 = 35
  47 = 10
  48 = 26
  49 = 24
  50 = 40
  51 = 21
  52 = 11
}

# Set header row (A1, B1) - bold/bordered style from the template is preserved
$ws.Range("A1").Value = $names[0]
$ws.Range("B1").Value = $names[1]

# Set cluster names (col A) and active case counts (col B) for rows 2..52
# in ascending alphabetical order of cluster name, matching the source data.
for ($r = 2; $r -le 52; $r++) {
  $ws.Cells.Item($r, 1).Value = $names[$r]
  $ws.Cells.Item($r, 2).Value = $values[$r]
}
